$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.197005748748779
$ws.Range("B1").Value = 4.473758220672607
$ws.Range("C1").Value = 7.138808250427246
$ws.Range("D1").Value = 7.250575542449951
$ws.Range("E1").Value = 5.306873798370361
